$wb = $excel.ActiveWorkbook

# --- About sheet: add EU/US biofuel cost adjustment note + conversion values ---
$about = $wb.Worksheets.Item("About")
$about.Range("A26").Value = "We adjust for the EU data by the ratio of EU:US pre-tax transportation biofuel costs (see file fuels/BFPaT for the EU and US models)."
$about.Range("A27").Value = "EU"
$about.Range("B27").Value = 0.000018152570386688024
$about.Range("A28").Value = "US"
$about.Range("B28").Value = 0.000012337034592036476

# --- ICtPSFfL sheet: scale the renewable-diesel incremental cost row by the EU:US ratio ---
$ws = $wb.Worksheets.Item("ICtPSFfL")
$ws.Range("B7:AK7").Formula = '=MAX(Calcs!B35,0)*(About!$B$27/About!$B$28)'
